$d = $word.ActiveDocument

# Locate the "{{修改内容}}" placeholder. It is a unique token in the
# document, spread across three runs: "{{", "修改内容", "}}". The middle
# run already carries the desired rFonts/hint formatting, so rather than
# collapsing all three runs into a brand new one (which would pick up the
# first run's formatting), we trim the "{{" and "}}" runs away and just
# retarget the middle run's text to "/".
$all = $d.Content
$found = $all.Find.Execute("{{修改内容}}", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $startPos = $all.Start
    $endPos = $all.End

    # Remove the trailing "}}" first so the earlier offsets stay valid.
    $closeBrace = $d.Range($endPos - 2, $endPos)
    $closeBrace.Delete()

    # Remove the leading "{{".
    $openBrace = $d.Range($startPos, $startPos + 2)
    $openBrace.Delete()

    # What's left between the (now gone) braces is the original
    # "修改内容" run; replace its text with "/", keeping its formatting.
    $middle = $d.Range($startPos, $startPos + 4)
    $middle.Text = "/"
}
